$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.578.90'
$ws.Range("E2").Value = '  +2.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.912.88'
$ws.Range("E3").Value = '  +5.54%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.67'
$ws.Range("E5").Value = '  +1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("E7").Value = '  +1.77%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3961'
$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09803'
$ws.Range("E9").Value = '  +1.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.164'
$ws.Range("E10").Value = '  +5.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.89'
$ws.Range("E11").Value = '  +2.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.552'
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.15'
$ws.Range("E13").Value = '  +3.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.908.00'
$ws.Range("E14").Value = '  +5.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.590'
$ws.Range("E15").Value = '  +4.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9996'
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001146'
$ws.Range("E17").Value = '  +1.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.77'
$ws.Range("E18").Value = '  +1.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06652'
$ws.Range("E19").Value = '  -0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.13'
$ws.Range("E20").Value = '  +5.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9996'
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.290'
$ws.Range("E22").Value = '  +6.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.624.73'
$ws.Range("E23").Value = '  +2.27%  '

$ws.Range("E24").Value = '  +2.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.282'
$ws.Range("E25").Value = '  +1.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.755'
$ws.Range("E26").Value = '  +15.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.130.01'
$ws.Range("E27").Value = '  +5.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.34'
$ws.Range("E28").Value = '  +3.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.48'
$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.97'
$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.108'
$ws.Range("E31").Value = '  +6.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1076'
$ws.Range("E32").Value = '  +1.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.732'
$ws.Range("E33").Value = '  +2.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.638'
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.857'
$ws.Range("E35").Value = '  +10.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06809'
$ws.Range("E36").Value = '  +1.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02450'
$ws.Range("E37").Value = '  +5.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.276'
$ws.Range("E38").Value = '  +10.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2238'
$ws.Range("E39").Value = '  +4.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.108'

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.69'
$ws.Range("E41").Value = '  +3.99%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6468'
$ws.Range("E42").Value = '  +4.70%  '

$ws.Range("E43").Value = '  +4.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9992'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.75'
$ws.Range("E45").Value = '  +4.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6119'
$ws.Range("E46").Value = '  +3.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.290'
$ws.Range("E47").Value = '  +0.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.661'
$ws.Range("E48").Value = '  -0.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.057'
$ws.Range("E49").Value = '  +6.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.47'
$ws.Range("E50").Value = '  +1.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.213'
$ws.Range("E51").Value = '  +2.92%  '
